# E-commerce.xlsx workbook edit: "adding models and routes"
#
# Rename the "Relaciones" worksheet to "Relaciones_SQL" and update the
# in-sheet selection that was left active on that sheet (H20 -> J19).

$wb = $excel.ActiveWorkbook

# Rename the "Relaciones" sheet to "Relaciones_SQL"
$ws = $wb.Worksheets.Item("Relaciones")
$ws.Name = "Relaciones_SQL"

# Make sure it's the active sheet and move the selection to J19,
# matching the saved cursor position recorded in the workbook.
$ws.Activate()
$ws.Range("J19").Select()
